$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.27"
$ws.Range("E2").Value = "'-0.19%"
$ws.Range("D3").Value = "'43.86"
$ws.Range("E3").Value = "'-0.18%"
$ws.Range("D4").Value = "'5.557"
$ws.Range("E4").Value = "'-0.42%"
$ws.Range("D5").Value = "'0.08055"
$ws.Range("E5").Value = "'-0.67%"
$ws.Range("D6").Value = "'1.972"
$ws.Range("E6").Value = "'3.73%"
$ws.Range("D7").Value = "'4.323"
$ws.Range("E7").Value = "'0.91%"
$ws.Range("D8").Value = "'2.553"
$ws.Range("D9").Value = "'0.9451"
$ws.Range("E9").Value = "'0.71%"
$ws.Range("D10").Value = "'0.1172"
$ws.Range("E10").Value = "'0.51%"
$ws.Range("D11").Value = "'0.1863"
$ws.Range("E11").Value = "'-1.83%"
$ws.Range("D12").Value = "'11.82"
$ws.Range("E12").Value = "'38.13%"
$ws.Range("D13").Value = "'0.09827"
$ws.Range("E13").Value = "'1.06%"
$ws.Range("D14").Value = "'0.04749"
$ws.Range("E14").Value = "'14.35%"
$ws.Range("E15").Value = "'-0.34%"
$ws.Range("D16").Value = "'0.001289"
$ws.Range("E16").Value = "'0.49%"
$ws.Range("D17").Value = "'0.04214"
$ws.Range("E17").Value = "'-2.31%"
$ws.Range("D18").Value = "'0.005936"
$ws.Range("E18").Value = "'-0.51%"
$ws.Range("E19").Value = "'-5.59%"
$ws.Range("D20").Value = "'0.3474"
$ws.Range("E20").Value = "'-0.31%"
$ws.Range("D21").Value = "'0.1408"
$ws.Range("E21").Value = "'3.34%"
$ws.Range("D22").Value = "'0.2508"
$ws.Range("E22").Value = "'-2.97%"
$ws.Range("D23").Value = "'0.001250"
$ws.Range("E23").Value = "'1.05%"
$ws.Range("D24").Value = "'0.004310"
$ws.Range("E24").Value = "'-2.08%"
$ws.Range("D25").Value = "'0.0001193"
$ws.Range("E25").Value = "'-3.06%"
$ws.Range("D26").Value = "'0.0003749"
$ws.Range("E26").Value = "'-5.98%"
$ws.Range("D38").Value = "'0.02584"
$ws.Range("E38").Value = "'-3.27%"
$ws.Range("D39").Value = "'0.05496"
$ws.Range("E39").Value = "'0.19%"
$ws.Range("D40").Value = "'0.007559"
$ws.Range("E40").Value = "'-1.87%"
$ws.Range("D41").Value = "'0.1401"
$ws.Range("E41").Value = "'0.05%"
$ws.Range("D42").Value = "'0.007463"
$ws.Range("E42").Value = "'-34.73%"
$ws.Range("D43").Value = "'0.002019"
$ws.Range("E43").Value = "'-4.19%"
$ws.Range("D44").Value = "'0.008352"
$ws.Range("E44").Value = "'-14.42%"
$ws.Range("D45").Value = "'0.00007089"
$ws.Range("E45").Value = "'1.15%"
$ws.Range("E46").Value = "'0.18%"
$ws.Range("E47").Value = "'1.44%"
$ws.Range("D48").Value = "'0.003522"
$ws.Range("E48").Value = "'1.60%"
$ws.Range("E49").Value = "'0.18%"
$ws.Range("E50").Value = "'0.18%"
